$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 12 data: A12 = 50000000, B12 = 1, C12 = A12/(B12*2)
$ws.Range("A12").Value = 50000000
$ws.Range("B12").Value = 1
$ws.Range("C12").Formula = "=A12/(B12*2)"

# Match the number formatting used by the rest of the factor column (C4:C11)
$ws.Range("C12").NumberFormat = $ws.Range("C11").NumberFormat

# Update the selected cell shown in the saved view
$ws.Range("I8").Select()
